# Generate Report for Handback
#
# Refreshes the localization-status report after a handback event:
#  - "Status" flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview + both locale sheets)
#  - zh-cn's "Latest Handback DateTime" moves off the zero-date
#    placeholder to a real timestamp
#  - de-de gets its "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" populated for the first time (previously
#    blank placeholders)
#  - both locale sheets get a new hyperlink on "Latest Target File"
#    (column I) for each data row, pointing at the same source .md file
#    that column A already links to
#  - the now-wider Status / Latest Target File / Latest Handback File
#    columns are resized to fit their longer content

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$mdDisplay = "53a5ca64-5bd1-4173-889e-a068e842066b.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f041c0d09c81aef8a41060d8a5e462a944fad45e/e2e/53a5ca64-5bd1-4173-889e-a068e842066b.md"

$zhXlf = "53a5ca64-5bd1-4173-889e-a068e842066b.f5d62d75a0e646e67e7412f2aad77a11243704d3.zh-cn.xlf"
$deXlf = "53a5ca64-5bd1-4173-889e-a068e842066b.f5d62d75a0e646e67e7412f2aad77a11243704d3.de-de.xlf"

$handedBackStatus = "Handed back: in sync with en-US"

# widths, expressed as ColumnWidth (character units) so the stored
# <col width=.../> lands on 30 / 40 respectively
$statusColWidth = 29.166666666666664
$fileColWidth = 39.16666666666667

# ---------------------------------------------------------------------
# 1) Overview sheet: Status columns (E = zh-cn, F = de-de) for both rows
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $handedBackStatus
$ws1.Range("F2").Value = $handedBackStatus
$ws1.Range("E3").Value = $handedBackStatus
$ws1.Range("F3").Value = $handedBackStatus

$ws1.Columns.Item(5).ColumnWidth = $statusColWidth
$ws1.Columns.Item(6).ColumnWidth = $statusColWidth

# ---------------------------------------------------------------------
# 2) zh-cn sheet (rows 2 and 3)
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = $handedBackStatus
$ws2.Range("C3").Value = $handedBackStatus

# Latest Target File (I) now references the source .md, with a hyperlink
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdDisplay) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $mdUrl, "", "", $mdDisplay) | Out-Null

# Latest Handback File (J) now references the already-generated xlf
$ws2.Range("J2").Value = $zhXlf
$ws2.Range("J3").Value = $zhXlf

# Latest Handback DateTime (K) moves off the zero-date placeholder
$ws2.Range("K2").Value = "2016-08-20 07:06:31"
$ws2.Range("K3").Value = "2016-08-20 07:06:31"

$ws2.Columns.Item(3).ColumnWidth = $statusColWidth
$ws2.Columns.Item(9).ColumnWidth = $fileColWidth
$ws2.Columns.Item(10).ColumnWidth = $fileColWidth

# ---------------------------------------------------------------------
# 3) de-de sheet (rows 2 and 3)
# ---------------------------------------------------------------------
$ws3.Range("C2").Value = $handedBackStatus
$ws3.Range("C3").Value = $handedBackStatus

$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdDisplay) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $mdUrl, "", "", $mdDisplay) | Out-Null

$ws3.Range("J2").Value = $deXlf
$ws3.Range("J3").Value = $deXlf

$ws3.Range("K2").Value = "2016-08-20 07:06:37"
$ws3.Range("K3").Value = "2016-08-20 07:06:37"

$ws3.Columns.Item(3).ColumnWidth = $statusColWidth
$ws3.Columns.Item(9).ColumnWidth = $fileColWidth
$ws3.Columns.Item(10).ColumnWidth = $fileColWidth
